$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" ---
# This status string appears on all three sheets:
#   Overview!E2, Overview!F2 (zh-cn / de-de status columns)
#   zh-cn!C2, de-de!C2 (Status column of each per-locale table)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Column width changes (status columns shrink to fit the shorter text) ---
# Target stored width ~= 13.4101845877511 character-units; the closest width
# this COM surface can produce (snapped to whole-pixel granularity) is
# 13.333333333333334, which is what ColumnWidth = 12.5 yields.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
